$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 542, shifting existing rows 542:628 down to 543:629
$ws.Rows.Item(542).Insert()

# Fill in the new row 542 with the new weekly record.
# Columns that stay identical across this whole market/product block
# (A,B,C,E,F,G,H,N,Q,R) are copied from the row directly below (row 543,
# which used to be row 542 before the insert).
$ws.Cells.Item(542, 1).Value2 = $ws.Cells.Item(543, 1).Value2   # A Mercado ID
$ws.Cells.Item(542, 2).Value2 = $ws.Cells.Item(543, 2).Value2   # B Mercado
$ws.Cells.Item(542, 3).Value2 = $ws.Cells.Item(543, 3).Value2   # C Region
$ws.Cells.Item(542, 4).Value2 = 44951                           # D Fecha
$ws.Cells.Item(542, 5).Value2 = $ws.Cells.Item(543, 5).Value2   # E Codreg
$ws.Cells.Item(542, 6).Value2 = $ws.Cells.Item(543, 6).Value2   # F Categoria ID
$ws.Cells.Item(542, 7).Value2 = $ws.Cells.Item(543, 7).Value2   # G Categoria
$ws.Cells.Item(542, 8).Value2 = $ws.Cells.Item(543, 8).Value2   # H Variedad
$ws.Cells.Item(542, 9).Value2 = "Primera"                       # I Calidad
$ws.Cells.Item(542, 10).Value2 = 250                            # J Volumen
$ws.Cells.Item(542, 11).Value2 = 3000                           # K Precio minimo
$ws.Cells.Item(542, 12).Value2 = 3300                           # L Precio maximo
$ws.Cells.Item(542, 13).Value2 = 3144                           # M Precio promedio ponderado
$ws.Cells.Item(542, 14).Value2 = $ws.Cells.Item(543, 14).Value2 # N Unidad de comercializacion
$ws.Cells.Item(542, 15).Value2 = $ws.Cells.Item(543, 15).Value2 # O Origen
$ws.Cells.Item(542, 16).Value2 = 87                             # P Precio $/Kg
$ws.Cells.Item(542, 17).Value2 = $ws.Cells.Item(543, 17).Value2 # Q Kg o Unidades
$ws.Cells.Item(542, 18).Value2 = $ws.Cells.Item(543, 18).Value2 # R Clasificacion

# Preserve the date style used by the rest of column D
$ws.Cells.Item(542, 4).NumberFormat = $ws.Cells.Item(543, 4).NumberFormat
